$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 45
$ws.Range("G45").Value = 1.42
$ws.Range("H45").Value = 4.1
$ws.Range("I45").Value = 8.5
$ws.Range("J45").Value = 2
$ws.Range("L45").Value = 7.5
$ws.Range("U45").Value = 2.25
$ws.Range("V45").Value = 1.57
$ws.Range("W45").Value = 5.5
$ws.Range("Y45").Value = 8.5
$ws.Range("Z45").Value = 9
$ws.Range("AE45").Value = 23
$ws.Range("AJ45").Value = 26
$ws.Range("AK45").Value = 101
$ws.Range("AL45").Value = 67
$ws.Range("AN45").Value = 3.2
$ws.Range("AO45").Value = 7
$ws.Range("AQ45").Value = 21
$ws.Range("AS45").Value = 201
$ws.Range("AV45").Value = 81
$ws.Range("AX45").Value = 8.5
$ws.Range("BA45").Value = 201

# Row 46
$ws.Range("M46").Value = 1.07
$ws.Range("N46").Value = 9
$ws.Range("Q46").Value = 2.25
$ws.Range("R46").Value = 1.62

# Row 47
$ws.Range("G47").Value = 1.95
$ws.Range("I47").Value = 4.5
$ws.Range("J47").Value = 2.75
$ws.Range("M47").Value = 1.13
$ws.Range("N47").Value = 6
$ws.Range("U47").Value = 2.38
$ws.Range("V47").Value = 1.53
$ws.Range("W47").Value = 5
$ws.Range("X47").Value = 7.5
$ws.Range("AI47").Value = 21
$ws.Range("AU47").Value = 10
$ws.Range("AX47").Value = 6

# Row 48
$ws.Range("K48").Value = 2.25
$ws.Range("L48").Value = 5.5
$ws.Range("M48").Value = 1.06
$ws.Range("N48").Value = 10
$ws.Range("O48").Value = 1.29
$ws.Range("P48").Value = 3.5
$ws.Range("Q48").Value = 1.93
$ws.Range("R48").Value = 1.93
$ws.Range("S48").Value = 1.36
$ws.Range("T48").Value = 3
$ws.Range("U48").Value = 1.83
$ws.Range("V48").Value = 1.83
$ws.Range("W48").Value = 7
$ws.Range("X48").Value = 7.5
$ws.Range("Z48").Value = 12
$ws.Range("AA48").Value = 13
$ws.Range("AB48").Value = 26
$ws.Range("AC48").Value = 10
$ws.Range("AE48").Value = 17
$ws.Range("AF48").Value = 51
$ws.Range("AG48").Value = 301
$ws.Range("AH48").Value = 15
$ws.Range("AJ48").Value = 19
$ws.Range("AL48").Value = 41
$ws.Range("AN48").Value = 3.5
$ws.Range("AP48").Value = 19
$ws.Range("AS48").Value = 151
$ws.Range("AT48").Value = 3
$ws.Range("AU48").Value = 8.5
$ws.Range("AV48").Value = 51
$ws.Range("AX48").Value = 7
$ws.Range("AY48").Value = 29
$ws.Range("AZ48").Value = 34
$ws.Range("BA48").Value = 101
$ws.Range("BB48").Value = 126
$ws.Range("BC48").Value = 251

# Row 62
$ws.Range("G62").Value = 1.8
$ws.Range("H62").Value = 3.4
$ws.Range("I62").Value = 4.33
$ws.Range("J62").Value = 2.5
$ws.Range("K62").Value = 2.1
$ws.Range("L62").Value = 4.75
$ws.Range("M62").Value = 1.07
$ws.Range("N62").Value = 9
$ws.Range("X62").Value = 8.5
$ws.Range("AC62").Value = 9
$ws.Range("AG62").Value = 301
$ws.Range("AH62").Value = 12
$ws.Range("AK62").Value = 51
$ws.Range("AL62").Value = 41
$ws.Range("AN62").Value = 3.75
$ws.Range("AO62").Value = 9.5
$ws.Range("AU62").Value = 8.5
$ws.Range("AY62").Value = 23
$ws.Range("AZ62").Value = 34
$ws.Range("BA62").Value = 81
$ws.Range("BC62").Value = 251

# Row 107
$ws.Range("Q107").Value = 2.1
$ws.Range("R107").Value = 1.7

# Row 108
$ws.Range("G108").Value = 2.75
$ws.Range("H108").Value = 2.9
$ws.Range("I108").Value = 2.63
$ws.Range("J108").Value = 3.75
$ws.Range("M108").Value = 1.13
$ws.Range("N108").Value = 6
$ws.Range("O108").Value = 1.53
$ws.Range("P108").Value = 2.38
$ws.Range("Q108").Value = 2.7
$ws.Range("R108").Value = 1.44
$ws.Range("S108").Value = 1.62
$ws.Range("T108").Value = 2.2
$ws.Range("U108").Value = 2.25
$ws.Range("V108").Value = 1.57
$ws.Range("Y108").Value = 12
$ws.Range("Z108").Value = 29
$ws.Range("AA108").Value = 29
$ws.Range("AC108").Value = 6
$ws.Range("AE108").Value = 21
$ws.Range("AF108").Value = 81
$ws.Range("AG108").Value = 1250
$ws.Range("AH108").Value = 6.5
$ws.Range("AI108").Value = 11
$ws.Range("AL108").Value = 29
$ws.Range("AN108").Value = 4.75
$ws.Range("AP108").Value = 34
$ws.Range("AQ108").Value = 67
$ws.Range("AR108").Value = 101
$ws.Range("AS108").Value = 351
$ws.Range("AT108").Value = 2.2
$ws.Range("AU108").Value = 9.5
$ws.Range("AV108").Value = 81
$ws.Range("AX108").Value = 4.5
$ws.Range("AZ108").Value = 34
$ws.Range("BC108").Value = 351

# Row 172
$ws.Range("G172").Value = 2.35
$ws.Range("I172").Value = 3.2
$ws.Range("J172").Value = 3.25
$ws.Range("L172").Value = 4
$ws.Range("M172").Value = 1.11
$ws.Range("N172").Value = 6.5
$ws.Range("W172").Value = 6
$ws.Range("X172").Value = 10
$ws.Range("Y172").Value = 10
$ws.Range("Z172").Value = 23
$ws.Range("AA172").Value = 23
$ws.Range("AE172").Value = 19
$ws.Range("AH172").Value = 7.5
$ws.Range("AI172").Value = 15
$ws.Range("AJ172").Value = 13
$ws.Range("AK172").Value = 34
$ws.Range("AL172").Value = 34
$ws.Range("AN172").Value = 4.33
$ws.Range("AX172").Value = 5
$ws.Range("AY172").Value = 19
$ws.Range("BA172").Value = 67

# Row 174
$ws.Range("G174").Value = 1.44
$ws.Range("H174").Value = 5
$ws.Range("I174").Value = 6
$ws.Range("N174").Value = 23
$ws.Range("Q174").Value = 1.33
$ws.Range("R174").Value = 3.4
$ws.Range("W174").Value = 13
$ws.Range("Y174").Value = 9
$ws.Range("AH174").Value = 26
$ws.Range("AI174").Value = 41
$ws.Range("AK174").Value = 67
$ws.Range("AQ174").Value = 17
$ws.Range("BB174").Value = 67
$ws.Range("BC174").Value = 101

# Row 175
$ws.Range("G175").Value = 2.72
$ws.Range("H175").Value = 3.4
$ws.Range("I175").Value = 2.3
$ws.Range("J175").Value = 3.3
$ws.Range("K175").Value = 2.18
$ws.Range("L175").Value = 2.87
$ws.Range("N175").Value = 7.9
$ws.Range("P175").Value = 3.6
$ws.Range("Q175").Value = 1.75
$ws.Range("R175").Value = 2
$ws.Range("S175").Value = 1.36
$ws.Range("T175").Value = 2.9
$ws.Range("V175").Value = 2.15
$ws.Range("W175").Value = 10.25
$ws.Range("X175").Value = 15
$ws.Range("Y175").Value = 10
$ws.Range("Z175").Value = 32
$ws.Range("AA175").Value = 21
$ws.Range("AB175").Value = 27
$ws.Range("AC175").Value = 7.9
$ws.Range("AD175").Value = 6.7
$ws.Range("AG175").Value = 350
$ws.Range("AI175").Value = 12.5
$ws.Range("AJ175").Value = 9
$ws.Range("AK175").Value = 24
$ws.Range("AL175").Value = 17.5
$ws.Range("AM175").Value = 25
$ws.Range("AN175").Value = 4.8
$ws.Range("AO175").Value = 14.5
$ws.Range("AP175").Value = 21
$ws.Range("AQ175").Value = 65
$ws.Range("AR175").Value = 90
$ws.Range("AS175").Value = 250
$ws.Range("AT175").Value = 2.9
$ws.Range("AU175").Value = 6.8
$ws.Range("AX175").Value = 4.4
$ws.Range("AY175").Value = 12
$ws.Range("AZ175").Value = 19
$ws.Range("BA175").Value = 45
$ws.Range("BB175").Value = 75
$ws.Range("BC175").Value = 200

# Row 177
$ws.Range("G177").Value = 2.18
$ws.Range("H177").Value = 2.9
$ws.Range("I177").Value = 3.5
$ws.Range("J177").Value = 2.72
$ws.Range("L177").Value = 4.05
$ws.Range("M177").Value = 1.06
$ws.Range("N177").Value = 8.77
$ws.Range("P177").Value = 2.65
$ws.Range("T177").Value = 2.45
$ws.Range("Y177").Value = 8.75
$ws.Range("AC177").Value = 7.5
$ws.Range("AD177").Value = 5.7
$ws.Range("AI177").Value = 18.5
$ws.Range("AJ177").Value = 12
$ws.Range("AK177").Value = 55
$ws.Range("AL177").Value = 35
$ws.Range("AO177").Value = 11.25
$ws.Range("AP177").Value = 18.5
$ws.Range("AV177").Value = 60
$ws.Range("AX177").Value = 5.3
$ws.Range("AY177").Value = 20
$ws.Range("BC177").Value = 350
